# Create a new "logical" worksheet after "comparison", populate it with the
# Elixir logical-operator reference table, and make it the active/selected
# sheet (mirrors the author adding a third tab to the operators workbook).

$wb = $excel.ActiveWorkbook

$comparisonSheet = $wb.Worksheets.Item("comparison")
$logicalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $comparisonSheet)
$logicalSheet.Name = "logical"

# Header row
$logicalSheet.Range("A1").Value = "Operator"
$logicalSheet.Range("B1").Value = "Description"

# Column A first (operators), then column B (descriptions) so the shared
# string table is populated in the same order the original author typed it.
$logicalSheet.Range("A2").Value = "and"
$logicalSheet.Range("A3").Value = "or"
$logicalSheet.Range("A4").Value = "not"

$logicalSheet.Range("B2").Value = "Condition A and B are true"
$logicalSheet.Range("B3").Value = "Condition A or B are true"
$logicalSheet.Range("B4").Value = "Invert the boolean value"

# Leave the new sheet active, with B4 selected, and activate the tab so it
# becomes the workbook's visible/active sheet on open.
[void]$logicalSheet.Range("B4").Select()
[void]$logicalSheet.Activate()
